$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("E10").Value = "[Ludoff-Cont.Lóg.Prog CLP, Joel L.-Tec. Fundição, Ludoff-Comandos Eletricos, Tiago Freitas-M.S.R. ar Cond.]"
$ws.Range("F10").Value = "[Leonardo-Mec. Manut.Equip. ind., Paulo Rob.-Usin. CNC, Paulo Rob.-M.A.Comp.CAD / CAM, Leonardo-Retífica]"

# Row 11
$ws.Range("B11").Value = "Rogério-Elem. Máqui"
$ws.Range("C11").Value = "-"
$ws.Range("D11").Value = "[Paulo Rob.-M.A.Comp.CAD / CAM, Leonardo-Retífica, Paulo Rob.-Usin. CNC, Leonardo-Mec. Manut.Equip. ind.]"
$ws.Range("E11").Value = "[Aderci-Fresagem, Wellington-Trat. Termicos, Ivan-Tec. Soldagem, Valmir-Calderaria]"
$ws.Range("F11").Value = "[Gisele-Ens. Dest. Não Desti., Nilton-Metrologia 2, Humberto-C.pneumática, Aselmo-M. Motor Endot.]"

# Row 12
$ws.Range("B12").Value = "Rogério-Elem. Máqui"
$ws.Range("C12").Value = "[-, -, Tiago Freitas-M.S.R. ar Cond., -]"
$ws.Range("D12").Value = "[Ludoff-Comandos Eletricos, Leonardo-Retífica, Ludoff-Cont.Lóg.Prog CLP, Leonardo-Mec. Manut.Equip. ind.]"
$ws.Range("E12").Value = "[Aderci-Fresagem, Wellington-Trat. Termicos, Ivan-Tec. Soldagem, Valmir-Calderaria]"
$ws.Range("F12").Value = "[Gisele-Ens. Dest. Não Desti., Nilton-Metrologia 2, Humberto-C.pneumática, Aselmo-M. Motor Endot.]"

# Row 14
$ws.Range("C14").Value = "[-, -, Tiago Freitas-M.S.R. ar Cond., -]"
$ws.Range("D14").Value = "[Ludoff-Comandos Eletricos, Joel L.-Tec. Fundição, Ludoff-Cont.Lóg.Prog CLP, -]"
$ws.Range("E14").Value = "[Aderci-Fresagem, Wellington-Trat. Termicos, Ivan-Tec. Soldagem, Valmir-Calderaria]"
$ws.Range("F14").Value = "[Gisele-Ens. Dest. Não Desti., Nilton-Metrologia 2, Humberto-C.pneumática, Aselmo-M. Motor Endot.]"

# Row 15
$ws.Range("B15").Value = "-"
$ws.Range("C15").Value = "[-, -, Tiago Freitas-M.S.R. ar Cond., Joel L.-Tec. Fundição]"
$ws.Range("D15").Value = "[-, Joel L.-Tec. Fundição, -, -]"
$ws.Range("E15").Value = "[Aderci-Fresagem, Wellington-Trat. Termicos, Ivan-Tec. Soldagem, Valmir-Calderaria]"
$ws.Range("F15").Value = "[Gisele-Ens. Dest. Não Desti., Nilton-Metrologia 2, Humberto-C.pneumática, Paulo Rob.-Usin. CNC]"

# Row 16
$ws.Range("E16").Value = "[Paulo Rob.-M.A.Comp.CAD / CAM, Ludoff-Cont.Lóg.Prog CLP, Aselmo-M. Motor Endot., Ludoff-Comandos Eletricos]"
$ws.Range("F16").Value = "[Paulo Rob.-M.A.Comp.CAD / CAM, Leonardo-Mec. Manut.Equip. ind., Leonardo-Retífica, Paulo Rob.-Usin. CNC]"
